# Reorders the data rows (2-17) of the sheet according to the mapping
# derived from the diff. Each target row receives the full A:T row
# content (values) that used to live at the corresponding source row.
# This is a pure permutation of existing rows - no values are created
# or removed, only relocated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# target row -> source row (both refer to row numbers in the ORIGINAL sheet)
$mapping = @{
    2  = 16
    3  = 3
    4  = 4
    5  = 7
    6  = 6
    7  = 11
    8  = 12
    9  = 17
    10 = 2
    11 = 15
    12 = 9
    13 = 10
    14 = 14
    15 = 5
    16 = 8
    17 = 13
}

$firstCol = 1   # A
$lastCol  = 20  # T

# 1) Snapshot all source rows' values (A:T) before any writes happen,
#    so that overwriting a row used as a source for another target
#    doesn't corrupt the data.
$snapshot = @{}
foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    if (-not $snapshot.ContainsKey($sourceRow)) {
        $rowValues = @{}
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $rowValues[$c] = $ws.Cells.Item($sourceRow, $c).Value2
        }
        $snapshot[$sourceRow] = $rowValues
    }
}

# 2) Write the snapshot values into the target rows.
foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $rowValues = $snapshot[$sourceRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($targetRow, $c).Value = $rowValues[$c]
    }
}
